$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "318.40", "43.948.81").
# Force it to stay plain text so Excel does not silently coerce it to a
# number (which would also lose the trailing-zero / thousand-dot formatting),
# then restore the default "Normal" style so no stray formatting is left
# behind on cells that were not otherwise restyled.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '43.948.81'
$ws.Range('E2').Value = '  +2.32%  '
$ws.Range('D3').Value = '2.254.42'
$ws.Range('E3').Value = '  +1.48%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '318.40'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').Value = '101.56'
$ws.Range('E6').Value = '  +2.39%  '
$ws.Range('E7').Value = '  -0.71%  '
$ws.Range('E9').Value = '  -0.78%  '
$ws.Range('D10').Value = '37.19'
$ws.Range('E10').Value = '  +1.08%  '
$ws.Range('D11').Value = '0.0830'
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('D12').Value = '7.62'
$ws.Range('E12').Value = '  -0.25%  '
$ws.Range('E13').Value = '  -1.39%  '
$ws.Range('D14').Value = '2.601.38'
$ws.Range('E14').Value = '  +1.36%  '
$ws.Range('D15').Value = '0.858'
$ws.Range('E15').Value = '  -0.32%  '
$ws.Range('D16').Value = '14.48'
$ws.Range('E16').Value = '  +1.04%  '
$ws.Range('D17').Value = '2.254.60'
$ws.Range('E17').Value = '  +1.45%  '
$ws.Range('D18').Value = '43.876.84'
$ws.Range('E18').Value = '  +2.29%  '
$ws.Range('D19').Value = '13.61'
$ws.Range('E19').Value = '  -5.93%  '
$ws.Range('D20').Value = '0.0₃0986'
$ws.Range('E20').Value = '  +2.39%  '
$ws.Range('D21').Value = '6.53'
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('D22').Value = '65.85'
$ws.Range('E22').Value = '  +1.24%  '
$ws.Range('D23').Value = '3.12'
$ws.Range('E23').Value = '  -1.70%  '
$ws.Range('D24').Value = '235.54'
$ws.Range('E24').Value = '  -0.45%  '
$ws.Range('D25').Value = '2.10'
$ws.Range('E25').Value = '  -3.32%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').Value = '10.20'
$ws.Range('E27').Value = '  +2.28%  '
$ws.Range('E28').Value = '  -3.49%  '
$ws.Range('D29').Value = '37.22'
$ws.Range('E29').Value = '  +3.92%  '
$ws.Range('D30').Value = '6.21'
$ws.Range('E30').Value = '  -2.11%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '20.17'
$ws.Range('E31').Value = '  -0.37%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').Value = '158.90'
$ws.Range('E32').Value = '  +3.13%  '
$ws.Range('D33').Value = '0.0852'
$ws.Range('E33').Value = '  -1.91%  '
$ws.Range('D34').Value = '2.69'
$ws.Range('E34').Value = '  +2.68%  '
$ws.Range('D35').Value = '0.116'
$ws.Range('E35').Value = '  +11.20%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').Value = '3.07'
$ws.Range('E36').Value = '  -1.83%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = '1.95'
$ws.Range('E37').Value = '  +1.78%  '
$ws.Range('E38').Value = '  -2.10%  '
$ws.Range('D39').Value = '16.33'
$ws.Range('E39').Value = '  +19.69%  '
$ws.Range('D40').Value = '3.72'
$ws.Range('E40').Value = '  +2.67%  '
$ws.Range('D41').Value = '4.20'
$ws.Range('E41').Value = '  -5.10%  '
$ws.Range('E42').Value = '  -2.26%  '
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').Value = '1.808.45'
$ws.Range('E44').Value = '  +3.78%  '
$ws.Range('D45').Value = '75.71'
$ws.Range('E45').Value = '  +0.76%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').Value = '0.198'
$ws.Range('E46').Value = '  -2.63%  '
$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D47').Value = '82.44'
$ws.Range('E47').Value = '  -2.85%  '
$ws.Range('E48').Value = '  -1.06%  '
$ws.Range('D49').Value = '104.38'
$ws.Range('E49').Value = '  +1.34%  '
$ws.Range('D50').Value = '1.69'
$ws.Range('E50').Value = '  +8.04%  '
$ws.Range('D51').Value = '58.23'
$ws.Range('E51').Value = '  -0.17%  '

$priceRange.Style = "Normal"
